$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test items for "Nytt bilmärke" (New Car Brand) / "Ny bilmodell" (New Car Model)
$ws.Range("B4").Value = "Ja"
$ws.Range("B5").Value = "Ja"

$ws.Range("D6").Value = "New Car Brand knappen"
$ws.Range("D5").Value = "Se till att man kan skriva in i nytt land rutan och att knapparna är synliga när de ska"
$ws.Range("E5").Value = "ja"
$ws.Range("E6").Value = "ja"
$ws.Range("D7").Value = "New Car Model fönster dyker upp"
$ws.Range("D8").Value = "Man kan skriva in en ny bilmodell"

# Column D width to match design (closest achievable given engine's width quantization)
$ws.Columns.Item(4).ColumnWidth = 74.45

# Update selection to match the new active cell
$ws.Range("D8").Select()
